# Fix the "Sebastien Beague" bullet on the "Cohesive Teamwork" slide:
#   - correct the misspelled surname "Beague" -> "Beauge"
#   - expand the task description to mention "Social media Research"
#
# Net result for the paragraph:
#   "Sebastien Beague: Work on Contact Page"
#   -> "Sebastien Beauge: Social media Research and Work on Contact Page"

$p = $ppt.ActivePresentation

# Locate the slide / shape / paragraph that holds the "Sebastien Beague" bullet
# instead of hard-coding indices, so the script is resilient to minor reordering.
$targetSlide = $null
$targetShape = $null
$targetParaIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi, 1)
                if ($para.Text -like "*Beague*") {
                    $targetSlide = $slide
                    $targetShape = $shape
                    $targetParaIndex = $pi
                }
            }
        }
    }
}

if ($targetShape -ne $null -and $targetParaIndex -gt 0) {
    $tr = $targetShape.TextFrame.TextRange
    $para = $tr.Paragraphs($targetParaIndex, 1)

    # The paragraph currently consists of three runs:
    #   1) "Sebastien "
    #   2) "Beague"          (misspelled, flagged err="1")
    #   3) ": Work on Contact Page"
    $run2 = $para.Runs(2, 1)

    # Remove run 2 ("Beague") entirely by deleting its exact character span -
    # this collapses the now-empty run out of the paragraph.
    $toDelete = $tr.Characters($run2.Start, $run2.Text.Length)
    $toDelete.Text = ""

    # Re-fetch the paragraph (now down to two runs) and fix run 1's text so the
    # name reads "Sebastien Beauge".
    $para = $tr.Paragraphs($targetParaIndex, 1)
    $run1 = $para.Runs(1, 1)
    $run1.Text = "Sebastien Beauge"

    # Update run 2 (previously run 3) to add the "Social media Research and" text.
    $para = $tr.Paragraphs($targetParaIndex, 1)
    $run2 = $para.Runs(2, 1)
    $run2.Text = ": Social media Research and Work on Contact Page"
}
